# Weekly price-sheet update: a new record for the most recent week is
# inserted at the top of the data block (row 172), pushing the existing
# rows 172-174 down to 173-175.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 172.
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row 172 with the new weekly entry.
$ws.Cells.Item(172, 1).Value = 4
$ws.Cells.Item(172, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(172, 3).Value = "Los Lagos"
$ws.Cells.Item(172, 4).Value = 44628
$ws.Cells.Item(172, 5).Value = 10
$ws.Cells.Item(172, 6).Value = "Fruta"
$ws.Cells.Item(172, 7).Value = 100108
$ws.Cells.Item(172, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(172, 9).Value = 100108002
$ws.Cells.Item(172, 10).Value = "Mango"
$ws.Cells.Item(172, 11).Value = "Sin especificar"
$ws.Cells.Item(172, 12).Value = "Primera"
$ws.Cells.Item(172, 13).Value = 240
$ws.Cells.Item(172, 14).Value = 8000
$ws.Cells.Item(172, 15).Value = 8000
$ws.Cells.Item(172, 16).Value = 8000
$ws.Cells.Item(172, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(172, 18).Value = "Perú"
$ws.Cells.Item(172, 19).Value = 2000
$ws.Cells.Item(172, 20).Value = 4
